# Insert one new data row right before the "2026/12/29" block (old row 668)
# so that it becomes the new row 668, shifting all following rows down by
# one (old row 668 -> new row 669, ..., old row 709 -> new row 710).
# Then populate the newly inserted row with the new data point:
#   2026/01/19, 月, 7, 201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 668..709 down to 669..710, leaving row 668 blank and ready
$ws.Rows(668).Insert()

$dateCell = $ws.Cells.Item(668, 1)

# Column A holds dates formatted as plain text (e.g. "2026/01/19"), not real
# Excel date serials. Temporarily force a text number format so the
# assignment isn't auto-converted into a date value, then restore the
# "Normal" cell style so the cell ends up unstyled, matching its neighbours.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/19"
$dateCell.Style = "Normal"

$ws.Cells.Item(668, 2).Value = "月"
$ws.Cells.Item(668, 3).Value = 7
$ws.Cells.Item(668, 4).Value = 201

Write-Host "Inserted new row 668 (2026/01/19); rows 668-709 shifted down to 669-710."
